$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.966.97'
$ws.Range("E2").Value = '  -0.14%  '

# Row 3
$ws.Range("D3").Value = '2.116.62'
$ws.Range("E3").Value = '  +0.67%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +0.75%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '346.83'
$ws.Range("E5").Value = '  +0.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.66%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5199'
$ws.Range("E7").Value = '  +0.86%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4438'
$ws.Range("E8").Value = '  +0.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.78'
$ws.Range("E9").Value = '  +2.42%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09330'
$ws.Range("E10").Value = '  -0.61%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.182'
$ws.Range("E11").Value = '  +0.53%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.15'
$ws.Range("E12").Value = '  -0.59%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.509'
$ws.Range("E13").Value = '  +3.94%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.900'
$ws.Range("E14").Value = '  +2.20%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.086.24'
$ws.Range("E15").Value = '  -0.53%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '103.16'
$ws.Range("E16").Value = '  +3.51%  '

# Row 17
$ws.Range("E17").Value = '  -0.04%  '

# Row 18
$ws.Range("E18").Value = '  +0.64%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.49'
$ws.Range("E19").Value = '  +3.85%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06704'
$ws.Range("E20").Value = '  +0.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.304'
$ws.Range("E21").Value = '  +1.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.009'
$ws.Range("E22").Value = '  +0.55%  '

# Row 23
$ws.Range("D23").Value = '30.007.62'
$ws.Range("E23").Value = '  -0.28%  '

# Row 24
$ws.Range("E24").Value = '  +0.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.330'
$ws.Range("E25").Value = '  -0.12%  '

# Row 26
$ws.Range("D26").Value = '2.382.85'
$ws.Range("E26").Value = '  +1.65%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.10'
$ws.Range("E27").Value = '  +0.15%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.538'
$ws.Range("E28").Value = '  -0.71%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.61'
$ws.Range("E29").Value = '  -0.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.26'
$ws.Range("E30").Value = '  +0.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.146'
$ws.Range("E31").Value = '  -2.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.765'
$ws.Range("E32").Value = '  +7.33%  '

# Row 33
$ws.Range("E33").Value = '  -0.26%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.249'
$ws.Range("E34").Value = '  +0.04%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.980'
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.569'
$ws.Range("E36").Value = '  +6.20%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.75'
$ws.Range("E37").Value = '  +5.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02625'
$ws.Range("E38").Value = '  +2.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06859'
$ws.Range("E39").Value = '  +1.24%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7078'
$ws.Range("E40").Value = '  +2.18%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.69'
$ws.Range("E41").Value = '  +0.99%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2245'
$ws.Range("E42").Value = '  -1.79%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.328'
$ws.Range("E43").Value = '  +1.55%  '

# Row 44
$ws.Range("E44").Value = '  +2.46%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.51'
$ws.Range("E45").Value = '  +1.95%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.357'
$ws.Range("E46").Value = '  +2.58%  '

# Row 47
$ws.Range("E47").Value = '  +0.54%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000358'
$ws.Range("E48").Value = '  +0.89%  '

# Row 49
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.633'
$ws.Range("E49").Value = '  -0.03%  '

# Row 50
$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.260'
$ws.Range("E50").Value = '  +8.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.225'
$ws.Range("E51").Value = '  +0.23%  '
